$wb = $excel.ActiveWorkbook

# This script applies the numeric corrections for the scheduled-runner
# market-data refresh across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Columns H..N = currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ).

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Animal Glue
$ws.Range("H5").Value = 151.5
$ws.Range("I5").Value = 19.428572
$ws.Range("K5").Value = 19.428572
$ws.Range("M5").Value = 95.571428

# Row 6: Antidote
$ws.Range("H6").Value = 399.1
$ws.Range("J6").Value = 682.6667
$ws.Range("L6").Value = 2048.0001
$ws.Range("N6").Value = -2272.0001

# Row 18: Growth Formula Beta
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = 0

# Row 28: Enchanted Silver Ink
$ws.Range("H28").Value = 873.75
$ws.Range("I28").Value = 891
$ws.Range("J28").Value = 822
$ws.Range("K28").Value = 891
$ws.Range("L28").Value = 822
$ws.Range("M28").Value = -406
$ws.Range("N28").Value = -1792

# Row 111: Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 15730.818
$ws.Range("I111").Value = 17974.75
$ws.Range("J111").Value = 9747
$ws.Range("K111").Value = 53924.25
$ws.Range("L111").Value = 29241
$ws.Range("M111").Value = -50857.25
$ws.Range("N111").Value = -35375

# Row 116: Growth Formula Kappa
$ws.Range("H116").Value = 7644.4165
$ws.Range("I116").Value = 9324.666999999999
$ws.Range("K116").Value = 9324.666999999999
$ws.Range("M116").Value = -5882.666999999999

# Row 138: Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3346.4211
$ws.Range("J138").Value = 3575.923
$ws.Range("L138").Value = 10727.769
$ws.Range("N138").Value = -21007.769

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Bronze Plate
$ws.Range("H4").Value = 474.5
$ws.Range("I4").Value = 474.5
$ws.Range("K4").Value = 474.5
$ws.Range("M4").Value = -358.5

# Row 6: Bronze Hoplon
$ws.Range("H6").Value = 2171
$ws.Range("I6").Value = 2171
$ws.Range("K6").Value = 2171
$ws.Range("M6").Value = -1998

# Row 18: Brass Alembic
$ws.Range("H18").Value = 14999
$ws.Range("I18").Value = 14999
$ws.Range("K18").Value = 14999
$ws.Range("M18").Value = -14677

# Row 37: Steel Chainmail
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("N37").Value = 0

# Row 45: Mythril Ingot
$ws.Range("H45").Value = 3936.5
$ws.Range("I45").Value = 4356
$ws.Range("K45").Value = 4356
$ws.Range("M45").Value = -3979

# Row 80: Titanium Hoplon
$ws.Range("H80").Value = 132043.8
$ws.Range("I80").Value = 70109
$ws.Range("J80").Value = 147527.5
$ws.Range("K80").Value = 70109
$ws.Range("L80").Value = 147527.5
$ws.Range("M80").Value = -69111
$ws.Range("N80").Value = -149523.5

# Row 83: Titanium Hoplon
$ws.Range("H83").Value = 132043.8
$ws.Range("I83").Value = 70109
$ws.Range("J83").Value = 147527.5
$ws.Range("K83").Value = 210327
$ws.Range("L83").Value = 442582.5
$ws.Range("M83").Value = -205335
$ws.Range("N83").Value = -452566.5

$ws = $wb.Worksheets.Item("BSM")
# Row 43: Steel Scythe
$ws.Range("H43").Value = 495000
$ws.Range("J43").Value = 495000
$ws.Range("L43").Value = 495000
$ws.Range("N43").Value = -495362

# Row 107: Deepgold Nugget
$ws.Range("H107").Value = 3244.8333
$ws.Range("I107").Value = 3244.8333
$ws.Range("K107").Value = 3244.8333
$ws.Range("M107").Value = -1324.8333

# Row 134: Ruthenium Ingot
$ws.Range("H134").Value = 1645
$ws.Range("I134").Value = 1645
$ws.Range("K134").Value = 4935
$ws.Range("M134").Value = -2400

$ws = $wb.Worksheets.Item("CRP")
# Row 86: Birch Lumber
$ws.Range("H86").Value = 7340.636
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89: Birch Lumber
$ws.Range("H89").Value = 7340.636
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 107: White Oak Lumber
$ws.Range("H107").Value = 729.3
$ws.Range("J107").Value = 1047.5
$ws.Range("L107").Value = 1047.5
$ws.Range("N107").Value = -4887.5

$ws = $wb.Worksheets.Item("CUL")
# Row 38: Dark Vinegar
$ws.Range("H38").Value = 37.666668
$ws.Range("I38").Value = 41.25
$ws.Range("J38").Value = 9
$ws.Range("K38").Value = 123.75
$ws.Range("L38").Value = 27
$ws.Range("M38").Value = 223.25
$ws.Range("N38").Value = -721

# Row 41: Cornbread
$ws.Range("H41").Value = 2099
$ws.Range("I41").Value = 2099
$ws.Range("K41").Value = 6297
$ws.Range("M41").Value = -5959

# Row 64: Baked Onion Soup
$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Row 67: Baked Onion Soup
$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# Row 68: Fermented Butter
$ws.Range("H68").Value = 518
$ws.Range("J68").Value = 518
$ws.Range("L68").Value = 1554
$ws.Range("N68").Value = -3176

# Row 71: Fermented Butter
$ws.Range("H71").Value = 518
$ws.Range("J71").Value = 518
$ws.Range("L71").Value = 4662
$ws.Range("N71").Value = -12774

# Row 86: Birch Syrup
$ws.Range("H86").Value = 188
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89: Birch Syrup
$ws.Range("H89").Value = 188
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 107: Frantoio Oil
$ws.Range("H107").Value = 520.5
$ws.Range("J107").Value = 520.5
$ws.Range("L107").Value = 1561.5
$ws.Range("N107").Value = -5401.5

# Row 113: Night Vinegar
$ws.Range("H113").Value = 1787.8
$ws.Range("I113").Value = 650
$ws.Range("J113").Value = 2072.25
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 6216.75
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -10556.75

# Row 131: Tsai tou Vounou
$ws.Range("H131").Value = 528547.1
$ws.Range("I131").Value = 1339.1428
$ws.Range("K131").Value = 4017.4284
$ws.Range("M131").Value = 1022.5716

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper Ingot
$ws.Range("H2").Value = 93.75
$ws.Range("J2").Value = 350
$ws.Range("L2").Value = 350
$ws.Range("N2").Value = -576

# Row 18: Brass Gorget
$ws.Range("H18").Value = 14005
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 43: Malachite Earrings
$ws.Range("H43").Value = 1679.8889
$ws.Range("J43").Value = 14999
$ws.Range("L43").Value = 14999
$ws.Range("N43").Value = -15301

# Row 70: Mythrite Ingot
$ws.Range("H70").Value = 4999.5
$ws.Range("I70").Value = 4999.5
$ws.Range("K70").Value = 4999.5
$ws.Range("M70").Value = -4729.5

# Row 73: Mythrite Ingot
$ws.Range("H73").Value = 4999.5
$ws.Range("I73").Value = 4999.5
$ws.Range("K73").Value = 4999.5
$ws.Range("M73").Value = -4063.5

# Row 80: Hardsilver Ingot
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3000
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = 0

# Row 83: Hardsilver Ingot
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15000
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = 0

# Row 107: Hard Mudstone Whetstone
$ws.Range("H107").Value = 1349.875
$ws.Range("I107").Value = 531.7273
$ws.Range("K107").Value = 531.7273
$ws.Range("M107").Value = 1388.2727

# Row 132: Lar Ingot
$ws.Range("H132").Value = 10080.706
$ws.Range("I132").Value = 10126.6
$ws.Range("J132").Value = 9736.5
$ws.Range("K132").Value = 30379.8
$ws.Range("L132").Value = 29209.5
$ws.Range("M132").Value = -27849.8
$ws.Range("N132").Value = -34269.5

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Leather Calot
$ws.Range("H2").Value = 37128.145
$ws.Range("I2").Value = 45979.6
$ws.Range("J2").Value = 14999.5
$ws.Range("K2").Value = 45979.6
$ws.Range("L2").Value = 14999.5
$ws.Range("M2").Value = -45867.6
$ws.Range("N2").Value = -15223.5

# Row 82: Dragon Leather
$ws.Range("H82").Value = 3499
$ws.Range("J82").Value = 3499
$ws.Range("L82").Value = 3499
$ws.Range("N82").Value = -4221

# Row 85: Dragon Leather
$ws.Range("H85").Value = 3499
$ws.Range("J85").Value = 3499
$ws.Range("L85").Value = 3499
$ws.Range("N85").Value = -5995

# Row 132: Silver Lobo Leather
$ws.Range("H132").Value = 3270.5454
$ws.Range("I132").Value = 2775.3333
$ws.Range("J132").Value = 5499
$ws.Range("K132").Value = 8325.999899999999
$ws.Range("L132").Value = 16497
$ws.Range("M132").Value = -5795.999899999999
$ws.Range("N132").Value = -21557

$ws = $wb.Worksheets.Item("WVR")
# Row 74: Ramie Robe of Casting
$ws.Range("H74").Value = 20902.625
$ws.Range("I74").Value = 18599.5
$ws.Range("K74").Value = 18599.5
$ws.Range("M74").Value = -17663.5

# Row 77: Ramie Robe of Casting
$ws.Range("H77").Value = 20902.625
$ws.Range("I77").Value = 18599.5
$ws.Range("K77").Value = 55798.5
$ws.Range("M77").Value = -51118.5

# Row 81: Crawler Silk
$ws.Range("H81").Value = 1254387.2
$ws.Range("J81").Value = 5000600.5
$ws.Range("L81").Value = 10001201
$ws.Range("N81").Value = -10003323

# Row 84: Crawler Silk
$ws.Range("H84").Value = 1254387.2
$ws.Range("J84").Value = 5000600.5
$ws.Range("L84").Value = 50006005
$ws.Range("N84").Value = -50016613

# Row 107: Bright Linen Yarn
$ws.Range("H107").Value = 1706.5416
$ws.Range("I107").Value = 1627.3
$ws.Range("J107").Value = 2102.75
$ws.Range("K107").Value = 4881.9
$ws.Range("L107").Value = 6308.25
$ws.Range("M107").Value = -2961.9
$ws.Range("N107").Value = -10148.25

# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 3945.25
$ws.Range("I132").Value = 4927
$ws.Range("K132").Value = 14781
$ws.Range("M132").Value = -12251

